# Update cryptocurrency price/volume data cells to reflect the latest
# GitHub Actions scrape (Thu Aug 15 05:46:08 UTC 2024), including the
# WhiteBITCoin / Hedera rank swap at rows 46-47.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.470.32"
$ws.Range("E2").Value = "  -3.58%  "
$ws.Range("D3").Value = "2.646.90"
$ws.Range("E3").Value = "  -2.07%  "
$ws.Range("E4").Value = "  -0.23%  "
$ws.Range("D5").Value = "'521.70"
$ws.Range("E5").Value = "  -0.93%  "
$ws.Range("D6").Value = "'144.33"
$ws.Range("E6").Value = "  -0.41%  "
$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").Value = "'0.569"
$ws.Range("E8").Value = "  -1.75%  "
$ws.Range("D9").Value = "'6.69"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("E10").Value = "  -3.00%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("E12").Value = "  +1.33%  "
$ws.Range("D13").Value = "3.113.57"
$ws.Range("E13").Value = "  -2.10%  "
$ws.Range("D14").Value = "58.475.74"
$ws.Range("E14").Value = "  -3.53%  "
$ws.Range("D15").Value = "'20.85"
$ws.Range("E15").Value = "  -2.14%  "
$ws.Range("D16").Value = "'0.0000136"
$ws.Range("E16").Value = "  -1.32%  "
$ws.Range("D17").Value = "2.656.86"
$ws.Range("E17").Value = "  -8.35%  "
$ws.Range("D18").Value = "'337.88"
$ws.Range("E18").Value = "  -2.95%  "
$ws.Range("D19").Value = "'4.41"
$ws.Range("E19").Value = "  -2.24%  "
$ws.Range("D20").Value = "'10.48"
$ws.Range("E20").Value = "  -1.35%  "
$ws.Range("E21").Value = "  -1.67%  "
$ws.Range("D23").Value = "'64.50"
$ws.Range("E23").Value = "  +1.35%  "
$ws.Range("D24").Value = "'0.424"
$ws.Range("E24").Value = "  +0.92%  "
$ws.Range("E25").Value = "  -1.81%  "
$ws.Range("D26").Value = "'0.999"
$ws.Range("E26").Value = "  +0.61%  "
$ws.Range("D27").Value = "0.0₃0797"
$ws.Range("E27").Value = "  -2.54%  "
$ws.Range("D28").Value = "'7.11"
$ws.Range("E28").Value = "  -2.69%  "
$ws.Range("D29").Value = "'6.63"
$ws.Range("E29").Value = "  -2.10%  "
$ws.Range("E30").Value = "  +0.02%  "
$ws.Range("E31").Value = "  -1.06%  "
$ws.Range("D32").Value = "'152.87"
$ws.Range("E32").Value = "  +1.69%  "
$ws.Range("D33").Value = "'18.86"
$ws.Range("E33").Value = "  -1.48%  "
$ws.Range("D34").Value = "'4.14"
$ws.Range("E34").Value = "  -2.22%  "
$ws.Range("E35").Value = "  -3.75%  "
$ws.Range("D36").Value = "'0.906"
$ws.Range("E36").Value = "  -3.92%  "
$ws.Range("D37").Value = "'0.856"
$ws.Range("E37").Value = "  -2.27%  "
$ws.Range("D38").Value = "'36.79"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("E39").Value = "  -4.96%  "
$ws.Range("E40").Value = "  -0.82%  "
$ws.Range("E41").Value = "  +0.47%  "
$ws.Range("D42").Value = "'0.607"
$ws.Range("E42").Value = "  -0.83%  "
$ws.Range("E43").Value = "  -2.15%  "
$ws.Range("D44").Value = "'269.83"
$ws.Range("E44").Value = "  -5.16%  "
$ws.Range("D45").Value = "'19.41"
$ws.Range("E45").Value = "  -3.78%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "'10.64"
$ws.Range("E46").Value = "  +1.59%  "
$ws.Range("B47").Value = "Hedera"
$ws.Range("C47").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D47").Value = "'0.0536"
$ws.Range("E47").Value = "  -0.54%  "
$ws.Range("D48").Value = "2.045.13"
$ws.Range("E48").Value = "  -4.48%  "
$ws.Range("E49").Value = "  -2.47%  "
$ws.Range("E50").Value = "  -3.06%  "
$ws.Range("D51").Value = "'18.36"
$ws.Range("E51").Value = "  -5.27%  "
